# Applies the "Updated symbol list" edit (Sat Feb  4 06:32:14 UTC 2023)
# to the cryptos worksheet: refreshed price/volume figures for the
# existing coin rows, plus a reshuffle where UpBots (previously a
# placeholder "--" row near the bottom) moves up to row 18 with live
# data, pushing LEO..BNIXToken down by one row, and BOLO /
# CoinbaseStockToken swap places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All of these cells hold plain text in the source sheet (t="inlineStr")
# even though many values look numeric/percentage/date-like. A leading
# apostrophe forces Excel to keep them as literal text instead of
# re-interpreting "329.20" -> 329.2, "1,779.16%" -> a number, etc.

$ws.Range("D2").Value = "'329.20"
$ws.Range("E2").Value = "'1.83%"
$ws.Range("D3").Value = "'41.20"
$ws.Range("E3").Value = "'3.92%"
$ws.Range("E4").Value = "'-4.40%"
$ws.Range("D5").Value = "'0.08166"
$ws.Range("E5").Value = "'1.75%"
$ws.Range("D6").Value = "'2.029"
$ws.Range("E6").Value = "'5.03%"
$ws.Range("D7").Value = "'8.774"
$ws.Range("E7").Value = "'1.19%"
$ws.Range("D8").Value = "'4.530"
$ws.Range("E8").Value = "'-0.98%"
$ws.Range("D9").Value = "'2.941"
$ws.Range("E9").Value = "'-0.23%"
$ws.Range("D10").Value = "'0.9192"
$ws.Range("E10").Value = "'-1.13%"
$ws.Range("D11").Value = "'0.1269"
$ws.Range("E11").Value = "'0.46%"
$ws.Range("D12").Value = "'0.1953"
$ws.Range("E12").Value = "'-0.75%"
$ws.Range("D13").Value = "'0.09255"
$ws.Range("E13").Value = "'1.36%"
$ws.Range("D14").Value = "'0.03739"
$ws.Range("E14").Value = "'3.56%"
$ws.Range("D15").Value = "'0.1060"
$ws.Range("E15").Value = "'1.31%"
$ws.Range("D16").Value = "'0.001303"
$ws.Range("E16").Value = "'0.53%"
$ws.Range("D17").Value = "'0.006184"
$ws.Range("E17").Value = "'-1.85%"
$ws.Range("B18").Value = "'UpBots"
$ws.Range("C18").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D18").Value = "'0.007503"
$ws.Range("E18").Value = "'1,779.16%"
$ws.Range("B19").Value = "'LEO"
$ws.Range("C19").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.437"
$ws.Range("E19").Value = "'2.64%"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3494"
$ws.Range("E20").Value = "'-1.19%"
$ws.Range("B21").Value = "'MCDex"
$ws.Range("C21").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'8.283"
$ws.Range("E21").Value = "'-4.97%"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "'0.1393"
$ws.Range("E22").Value = "'1.58%"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2391"
$ws.Range("E23").Value = "'-2.42%"
$ws.Range("B24").Value = "'CoinExToken"
$ws.Range("C24").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04415"
$ws.Range("E24").Value = "'0.07%"
$ws.Range("B25").Value = "'BitKan"
$ws.Range("C25").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").Value = "'0.001262"
$ws.Range("E25").Value = "'-0.05%"
$ws.Range("B26").Value = "'HotbitToken"
$ws.Range("C26").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").Value = "'0.004313"
$ws.Range("E26").Value = "'-2.03%"
$ws.Range("B27").Value = "'NitroEx"
$ws.Range("C27").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").Value = "'0.0001182"
$ws.Range("E27").Value = "'3.66%"
$ws.Range("B28").Value = "'Spectre.aiUtilityToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("B29").Value = "'LegolasExchange"
$ws.Range("C29").Value = "'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("B30").Value = "'BitZToken"
$ws.Range("C30").Value = "'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("B31").Value = "'Birake"
$ws.Range("C31").Value = "'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("B32").Value = "'NashExchange"
$ws.Range("C32").Value = "'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("B33").Value = "'AAXToken"
$ws.Range("C33").Value = "'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("B34").Value = "'CenX"
$ws.Range("C34").Value = "'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("B35").Value = "'BNIXToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("D39").Value = "'0.02756"
$ws.Range("E39").Value = "'9.09%"
$ws.Range("E40").Value = "'2.20%"
$ws.Range("D41").Value = "'0.007672"
$ws.Range("E41").Value = "'3.28%"
$ws.Range("D42").Value = "'0.1415"
$ws.Range("E42").Value = "'0.71%"
$ws.Range("D43").Value = "'0.009011"
$ws.Range("E43").Value = "'-6.20%"
$ws.Range("D44").Value = "'0.002124"
$ws.Range("E44").Value = "'0.31%"
$ws.Range("D45").Value = "'0.01135"
$ws.Range("E45").Value = "'13.69%"
$ws.Range("D46").Value = "'0.00006877"
$ws.Range("E46").Value = "'2.18%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002282"
$ws.Range("E48").Value = "'60.44%"
$ws.Range("B49").Value = "'BOLO"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003488"
$ws.Range("E49").Value = "'16.02%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.11%"
